# Update the "department" column (C) on the "courses" sheet:
# rows 2-4 (single-course rows) -> "Automotive"
# rows 5-6 (package rows) -> "Packages"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("courses")

$ws.Range("C2").Value = "Automotive"
$ws.Range("C3").Value = "Automotive"
$ws.Range("C4").Value = "Automotive"
$ws.Range("C5").Value = "Packages"
$ws.Range("C6").Value = "Packages"
